$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("C2").Style

$ws.Range('D2').Value = '30.545.79'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.881.13'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9995'
$ws.Range('D4').Style = $plainStyle
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.10'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4716'
$ws.Range('D7').Style = $plainStyle
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2883'
$ws.Range('D8').Style = $plainStyle
$ws.Range('E8').Value = '  -1.54%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06535'
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '22.06'
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '101.11'
$ws.Range('D11').Style = $plainStyle
$ws.Range('E11').Value = '  +4.17%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07819'
$ws.Range('D12').Style = $plainStyle
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.7492'
$ws.Range('D13').Style = $plainStyle
$ws.Range('E13').Value = '  +1.78%  '
$ws.Range('D14').Value = '1.877.09'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.223'
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '285.20'
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').Value = '30.545.58'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.18'
$ws.Range('D18').Style = $plainStyle
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007522'
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9995'
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '2.119.91'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.372'
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.9992'
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.390'
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.128'
$ws.Range('D25').Style = $plainStyle
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.81'
$ws.Range('D26').Style = $plainStyle
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.04'
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.913'
$ws.Range('D28').Style = $plainStyle
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.09698'
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.325'
$ws.Range('D30').Style = $plainStyle
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.494'
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.266'
$ws.Range('D32').Style = $plainStyle
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.190'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04838'
$ws.Range('D34').Style = $plainStyle
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.131'
$ws.Range('D35').Style = $plainStyle
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6947'
$ws.Range('D36').Style = $plainStyle
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.768'
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = '  +1.76%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01913'
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.866'
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '76.41'
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = '  +0.49%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.349'
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.981'
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4250'
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9989'
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8291'
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.43'
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.789'
$ws.Range('D47').Style = $plainStyle
$ws.Range('E47').Value = '  +2.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.048'
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '35.13'
$ws.Range('D49').Style = $plainStyle
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '889.50'
$ws.Range('D51').Style = $plainStyle
$ws.Range('E51').Value = '  -3.22%  '
